# Refresh the crypto price (column D) and 1h volume-change (column E)
# figures for rows 2-51, matching the latest source snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> @(new Price text, new Volume(1h) text, IsPriceNumericLooking)
$values = @{
    2 = @('30.566.50', '  +0.51%  ', $false)
    3 = @('2.138.29', '  +1.83%  ', $false)
    4 = @('1.009', '  +0.36%  ', $true)
    5 = @('352.02', '  +5.34%  ', $true)
    6 = @('1.007', '  +0.33%  ', $true)
    7 = @('0.5263', '  +0.97%  ', $true)
    8 = @('0.4564', '  +0.42%  ', $true)
    9 = @('53.70', '  -2.08%  ', $true)
    10 = @('0.09162', '  +3.15%  ', $true)
    11 = @('1.191', '  +1.11%  ', $true)
    12 = @('25.43', '  +5.75%  ', $true)
    13 = @('2.135.41', '  +1.66%  ', $false)
    14 = @('6.904', '  +1.70%  ', $true)
    15 = @('8.184', '  +2.23%  ', $true)
    16 = @('102.30', '  +5.41%  ', $true)
    17 = @('0.00001174', '  +2.57%  ', $true)
    18 = @('1.009', '  +0.39%  ', $true)
    19 = @('0.06734', '  +1.74%  ', $true)
    20 = @('20.58', '  +7.40%  ', $true)
    21 = @('1.007', '  +0.35%  ', $true)
    22 = @('6.385', '  +1.73%  ', $true)
    23 = @('30.659.48', '  +0.65%  ', $false)
    24 = @('12.94', '  +4.96%  ', $true)
    25 = @('2.381', '  +0.87%  ', $true)
    26 = @('2.391.68', '  +1.96%  ', $false)
    27 = @('22.54', '  +1.68%  ', $true)
    28 = @('2.645', '  +5.54%  ', $true)
    29 = @('165.01', '  +1.36%  ', $true)
    30 = @('136.33', '  +2.49%  ', $true)
    31 = @('1.225', '  +1.74%  ', $true)
    32 = @('0.1083', '  +1.51%  ', $true)
    33 = @('1.689', '  +2.62%  ', $true)
    34 = @('6.419', '  +0.86%  ', $true)
    35 = @('4.041', '  +2.50%  ', $true)
    36 = @('6.142', '  +5.98%  ', $true)
    37 = @('10.50', '  +1.25%  ', $true)
    38 = @('0.02651', '  +3.12%  ', $true)
    39 = @('0.06996', '  +2.33%  ', $true)
    40 = @('0.2341', '  +1.35%  ', $true)
    41 = @('12.80', '  +0.94%  ', $true)
    42 = @('0.7013', '  +2.28%  ', $true)
    43 = @('1.278', '  +2.54%  ', $true)
    44 = @('14.89', '  +7.07%  ', $true)
    45 = @('2.367', '  +2.36%  ', $true)
    46 = @('0.6538', '  +3.14%  ', $true)
    47 = @('0.00000000371', '  +8.89%  ', $true)
    48 = @('3.755', '  +2.88%  ', $true)
    49 = @('1.253', '  +0.53%  ', $true)
    50 = @('84.05', '  +1.36%  ', $true)
    51 = @('0.07300', '  +2.46%  ', $true)
}

foreach ($row in $values.Keys) {
    $priceText  = $values[$row][0]
    $volText    = $values[$row][1]
    $isNumeric  = $values[$row][2]

    $priceCell = $ws.Range("D$row")
    if ($isNumeric) {
        # Force text storage so a numeric-looking price (e.g. "1.009")
        # is not silently re-interpreted/reformatted as a number.
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $priceText
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $priceText
    }

    $ws.Range("E$row").Value = $volText
}
